$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New group (S17 / G03): Indicator-based alerts & holdings analytics
$groupDesc = "Indicator-based alerts & holdings analytics"

$startRow = 136

$tasks = @(
    @{
        Id = "S17_G03_TB001"
        Desc = "Design indicator_rule schema and alert condition model on top of candles store."
        Remarks = "IndicatorRule schema and JSON-based condition model defined on top of the candles store for indicator-driven alerts."
        Pending = "Extend schema for additional indicators (VWAP, MA cross) as needed."
    },
    @{
        Id = "S17_G03_TB002"
        Desc = "Add ORM models + Alembic migration for indicator_rules and alert/source fields."
        Remarks = "Alembic migration creates indicator_rules table and links alerts via rule_id and source=INTERNAL_INDICATOR."
        Pending = "Verify migrations on non-dev databases when rolling out."
    },
    @{
        Id = "S17_G03_TB003"
        Desc = "Implement indicator engine to compute RSI, MAs, volatility, ATR, and performance windows from candles."
        Remarks = "indicator_alerts service computes RSI, moving averages, volatility, ATR, performance windows, and volume ratios using market_data.load_series."
        Pending = "Add more indicators (e.g., VWAP, MA cross) and unit benchmarks if required."
    },
    @{
        Id = "S17_G03_TB004"
        Desc = "Implement rule evaluation + alert/order creation pipeline using existing risk engine and orders API."
        Remarks = "Rule evaluation engine resolves HOLDINGS universe, evaluates conditions with AND/OR logic, records INTERNAL_INDICATOR alerts, and enqueues optional WAITING orders (SELL_PERCENT / BUY_QUANTITY)."
        Pending = "Tighten risk-engine integration once live trading usage grows."
    },
    @{
        Id = "S17_G03_TB005"
        Desc = "Add scheduler / evaluation endpoint to run indicator rules periodically in IST market hours."
        Remarks = "Background scheduler thread evaluates enabled indicator rules every few minutes in IST; evaluate_indicator_rules_once() is available for manual runs."
        Pending = "Consider external scheduler/cron wiring for multi-process deployments."
    },
    @{
        Id = "S17_G03_TF001"
        Desc = "Expose indicator columns (RSI, MA%, volatility, performance) in Holdings DataGrid with filtering."
        Remarks = "Holdings DataGrid now derives RSI(14), 1M/1Y price performance, volatility 20D%, ATR(14)%, and volume-vs-20D-average columns from OHLCV history with numeric filters and negative-value highlighting."
        Pending = "Fine-tune default column visibility and add additional windows if desired."
    },
    @{
        Id = "S17_G03_TF002"
        Desc = "Add TradingView-style Alert modal on Holdings rows to create/edit indicator rules."
        Remarks = "Each Holdings row exposes an Alert button that opens a modal to configure indicator rules (indicator, operator, thresholds, period/window, trigger mode, and action type) backed by /api/indicator-alerts CRUD."
        Pending = "Extend modal to support multi-condition rules and per-rule naming."
    },
    @{
        Id = "S17_G03_TF003"
        Desc = "Surface fired indicator alerts and resulting WAITING orders in Alerts panel and Queue view."
        Remarks = "Fired indicator rules create INTERNAL_INDICATOR alerts and corresponding WAITING orders that flow through the existing queue and execution pipeline alongside TradingView alerts."
        Pending = "Add explicit UI badges/filters to distinguish indicator-based entries in the queue."
    }
)

for ($i = 0; $i -lt $tasks.Count; $i++) {
    $r = $startRow + $i
    $task = $tasks[$i]

    $ws.Cells.Item($r, 1).Value = "S17"
    $ws.Cells.Item($r, 2).Value = "G03"
    $ws.Cells.Item($r, 3).Value = $groupDesc
    $ws.Cells.Item($r, 4).Value = $task.Id
    $ws.Cells.Item($r, 5).Value = $task.Desc
    $ws.Cells.Item($r, 7).Value = "implemented"
    $ws.Cells.Item($r, 8).Value = $task.Remarks
    $ws.Cells.Item($r, 9).Value = $task.Pending
}

# Match the updated scroll position / active selection from the diff.
$ws.Range("C138").Select()
